$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and volume-change (E) figures, plus three pairs of
# rows (44/45, 48/49, 50/51) whose coin name/link/price/change were
# re-ordered by the latest coinranking.com pull.

# Cells whose new text is unambiguous (won't be auto-parsed as a number)
# -- plain assignment keeps them as text with no style change.
$plainTextUpdates = @(
    @{ Cell = "D2"; Value = "29.088.77" },
    @{ Cell = "E2"; Value = "  +0.29%  " },
    @{ Cell = "D3"; Value = "1.833.41" },
    @{ Cell = "E3"; Value = "  +0.14%  " },
    @{ Cell = "E4"; Value = "  +0.44%  " },
    @{ Cell = "E5"; Value = "  +1.01%  " },
    @{ Cell = "E6"; Value = "  +0.55%  " },
    @{ Cell = "E7"; Value = "  +0.37%  " },
    @{ Cell = "E8"; Value = "  -1.94%  " },
    @{ Cell = "E9"; Value = "  +0.63%  " },
    @{ Cell = "E10"; Value = "  +1.34%  " },
    @{ Cell = "E11"; Value = "  -0.24%  " },
    @{ Cell = "D12"; Value = "1.827.52" },
    @{ Cell = "E12"; Value = "  -0.19%  " },
    @{ Cell = "E13"; Value = "  +1.21%  " },
    @{ Cell = "E14"; Value = "  +0.80%  " },
    @{ Cell = "E15"; Value = "  +0.88%  " },
    @{ Cell = "E16"; Value = "  -3.53%  " },
    @{ Cell = "E17"; Value = "  +1.03%  " },
    @{ Cell = "D18"; Value = "29.063.75" },
    @{ Cell = "E18"; Value = "  +0.23%  " },
    @{ Cell = "E19"; Value = "  +2.41%  " },
    @{ Cell = "E20"; Value = "  +0.29%  " },
    @{ Cell = "E21"; Value = "  +0.41%  " },
    @{ Cell = "E22"; Value = "  -1.18%  " },
    @{ Cell = "E23"; Value = "  +0.31%  " },
    @{ Cell = "E24"; Value = "  +1.25%  " },
    @{ Cell = "E25"; Value = "  +2.35%  " },
    @{ Cell = "E26"; Value = "  +1.23%  " },
    @{ Cell = "E27"; Value = "  +0.47%  " },
    @{ Cell = "E28"; Value = "  +0.80%  " },
    @{ Cell = "E29"; Value = "  +2.00%  " },
    @{ Cell = "E30"; Value = "  +0.98%  " },
    @{ Cell = "E31"; Value = "  +5.33%  " },
    @{ Cell = "E32"; Value = "  +0.64%  " },
    @{ Cell = "E33"; Value = "  +1.52%  " },
    @{ Cell = "E34"; Value = "  +0.43%  " },
    @{ Cell = "E35"; Value = "  -0.57%  " },
    @{ Cell = "E36"; Value = "  -3.19%  " },
    @{ Cell = "D37"; Value = "1.228.71" },
    @{ Cell = "E37"; Value = "  -2.98%  " },
    @{ Cell = "E38"; Value = "  -0.06%  " },
    @{ Cell = "E39"; Value = "  -0.22%  " },
    @{ Cell = "E40"; Value = "  +6.68%  " },
    @{ Cell = "E41"; Value = "  -0.03%  " },
    @{ Cell = "E42"; Value = "  +0.35%  " },
    @{ Cell = "E43"; Value = "  +0.14%  " },
    @{ Cell = "B44"; Value = "Aave" },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave" },
    @{ Cell = "E44"; Value = "  +2.07%  " },
    @{ Cell = "B45"; Value = "BabyDogeCoin" },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge" },
    @{ Cell = "E45"; Value = "  -0.26%  " },
    @{ Cell = "E46"; Value = "  -0.11%  " },
    @{ Cell = "E47"; Value = "  +1.58%  " },
    @{ Cell = "B48"; Value = "XinFinNetwork" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc" },
    @{ Cell = "E48"; Value = "  +9.29%  " },
    @{ Cell = "B49"; Value = "EnergySwap" },
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" },
    @{ Cell = "E49"; Value = "  +0.98%  " },
    @{ Cell = "B50"; Value = "Cronos" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro" },
    @{ Cell = "E50"; Value = "  +1.01%  " },
    @{ Cell = "B51"; Value = "RenderToken" },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" },
    @{ Cell = "E51"; Value = "  +2.77%  " }
)

foreach ($u in $plainTextUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Cells whose new text LOOKS like a number (e.g. "1.003") -- force the
# Text format before writing, then restore the default "Normal" style so
# the cell keeps displaying/storing the literal string, with no residual
# number-format style difference versus the original cell.
$forcedTextUpdates = @(
    @{ Cell = "D4"; Value = "1.003" },
    @{ Cell = "D5"; Value = "243.62" },
    @{ Cell = "D6"; Value = "0.6285" },
    @{ Cell = "D8"; Value = "0.07466" },
    @{ Cell = "D9"; Value = "0.2929" },
    @{ Cell = "D10"; Value = "22.99" },
    @{ Cell = "D11"; Value = "0.07722" },
    @{ Cell = "D13"; Value = "5.008" },
    @{ Cell = "D14"; Value = "0.6674" },
    @{ Cell = "D15"; Value = "83.20" },
    @{ Cell = "D16"; Value = "0.000009408" },
    @{ Cell = "D17"; Value = "6.039" },
    @{ Cell = "D19"; Value = "12.60" },
    @{ Cell = "D20"; Value = "223.90" },
    @{ Cell = "D21"; Value = "1.004" },
    @{ Cell = "D22"; Value = "7.108" },
    @{ Cell = "D24"; Value = "160.06" },
    @{ Cell = "D25"; Value = "0.1401" },
    @{ Cell = "D26"; Value = "8.491" },
    @{ Cell = "D28"; Value = "1.499" },
    @{ Cell = "D29"; Value = "4.133" },
    @{ Cell = "D30"; Value = "4.060" },
    @{ Cell = "D31"; Value = "0.05449" },
    @{ Cell = "D32"; Value = "1.201" },
    @{ Cell = "D33"; Value = "0.7501" },
    @{ Cell = "D34"; Value = "1.849" },
    @{ Cell = "D35"; Value = "1.136" },
    @{ Cell = "D36"; Value = "2.609" },
    @{ Cell = "D38"; Value = "2.753" },
    @{ Cell = "D39"; Value = "0.01783" },
    @{ Cell = "D40"; Value = "6.637" },
    @{ Cell = "D41"; Value = "0.8942" },
    @{ Cell = "D43"; Value = "101.57" },
    @{ Cell = "D44"; Value = "65.52" },
    @{ Cell = "D45"; Value = "0.00000000123" },
    @{ Cell = "D46"; Value = "0.5101" },
    @{ Cell = "D47"; Value = "0.4039" },
    @{ Cell = "D48"; Value = "0.07435" },
    @{ Cell = "D49"; Value = "8.926" },
    @{ Cell = "D50"; Value = "0.05807" },
    @{ Cell = "D51"; Value = "1.656" }
)

foreach ($u in $forcedTextUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

